$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Remove the original "{c}" / " " row (row 3); everything below shifts up by one.
$ws.Rows("3:3").Delete()

# 2. Insert a new row for the "=" / "는" pair right before the "+" row
#    (which, after the deletion above, now sits at row 62), pushing
#    "+" and everything after it back down by one row.
$ws.Rows("62:62").Insert()
$ws.Cells.Item(62, 1).Value2 = "="
$ws.Cells.Item(62, 2).Value2 = "는"

# 3. Re-append the removed "{c}" / " " pair as a brand new row at the end (row 78).
$ws.Cells.Item(78, 1).Value2 = "{c}"
$ws.Cells.Item(78, 2).Value2 = " "

# Update the saved selection to match the author's final cursor position.
$ws.Range("B63").Select()
